# edit.ps1 - apply the report.docx changes described by the commit diff:
#   1. Merge the "TR" + "R" + "EB" hyperlink runs into a single "TRREB" run.
#   2. Merge the '\data\housing_market_Q1\raw' run with the following
#      closing-quote run into a single run.
#   3. Merge the '\data\housing_market_Q1\preprocessed' run with the
#      following closing-quote run into a single run, and add a new
#      level-2 list heading paragraph ("Extracting tables from the PDF")
#      right after that paragraph.

$d = $word.ActiveDocument
$rsq = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK ( ' )

# ---------------------------------------------------------------------
# 1) Hyperlink text "TR" + "R" + "EB" -> "TRREB"
#    Scope the Find/Replace strictly to the hyperlink's own Range so we
#    don't touch the many other "TR"/"TRREB" occurrences in the body text.
# ---------------------------------------------------------------------
$hyperlink = $d.Hyperlinks.Item(1)
$hlRange = $hyperlink.Range
$hlRange.Find.Execute("TRREB", $false, $false, $false, $false, $false, $true, 1, $false, "TRREB", 2)

# ---------------------------------------------------------------------
# 2) '\data\housing_market_Q1\raw' + closing quote -> single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "\data\housing_market_Q1\raw" + $rsq, $false, $false, $false, $false, $false,
    $true, 1, $false, "\data\housing_market_Q1\raw" + $rsq, 2)

# ---------------------------------------------------------------------
# 3) '\data\housing_market_Q1\preprocessed' + closing quote -> single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "\data\housing_market_Q1\preprocessed" + $rsq, $false, $false, $false, $false, $false,
    $true, 1, $false, "\data\housing_market_Q1\preprocessed" + $rsq, 2)

# ---------------------------------------------------------------------
# 4) Add a new paragraph right after the "...preprocessed'." paragraph
#    (the one that ends the "Modifying the PDF file" subsection), styled
#    like the other level-2 ("1.x") bold list headings, but additionally
#    justified ("both"), containing "Extracting tables from the PDF".
# ---------------------------------------------------------------------

# Locate the paragraph that now ends with the preprocessed-folder path
# sentence, by searching for its trailing text.
$anchor = $d.Content
$anchor.Find.Execute("\data\housing_market_Q1\preprocessed" + $rsq + ".", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchor.Paragraphs.First

# An existing level-2 heading paragraph ("Web scraping the PDF files.")
# supplies the numbering template (ilvl 1 / numId 1) we need to continue.
$templatePara = $d.Paragraphs.Item(10)

$newRange = $anchorPara.Range.InsertParagraphAfter()
$newParas = $d.Paragraphs
$newPara = $newParas.Item($anchorPara.Index + 1)

$newPara.Style = "List Paragraph"
$newPara.Alignment = 3
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($templatePara.Range.ListFormat.ListTemplate, $true, 2)
$newPara.Range.ListFormat.ListLevelNumber = 2
$newPara.Range.Bold = 1
$newPara.Range.BoldBi = 1
$newPara.Range.Text = "Extracting tables from the PDF"
